$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.027.15'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '2.754.22'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '352.00'
$ws.Range('E5').Value = '  -1.70%  '
$ws.Range('D6').Value = '107.18'
$ws.Range('E6').Value = '  -2.04%  '
$ws.Range('E7').Value = '  -2.68%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -1.94%  '
$ws.Range('D10').Value = '39.02'
$ws.Range('E10').Value = '  -2.51%  '
$ws.Range('E11').Value = '  +3.45%  '
$ws.Range('D12').Value = '0.0830'
$ws.Range('E12').Value = '  -2.12%  '
$ws.Range('D13').Value = '19.59'
$ws.Range('E13').Value = '  +0.54%  '
$ws.Range('D14').Value = '7.45'
$ws.Range('E14').Value = '  -2.20%  '
$ws.Range('D15').Value = '3.180.50'
$ws.Range('E15').Value = '  -0.88%  '
$ws.Range('D16').Value = '2.754.14'
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('D17').Value = '0.922'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').Value = '50.995.03'
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('D19').Value = '7.62'
$ws.Range('E19').Value = '  +3.23%  '
$ws.Range('D20').Value = '3.03'
$ws.Range('E20').Value = '  -2.60%  '
$ws.Range('D21').Value = '12.92'
$ws.Range('E21').Value = '  -1.34%  '
$ws.Range('D22').Value = '0.0₃0953'
$ws.Range('E22').Value = '  -2.57%  '
$ws.Range('D23').Value = '69.05'
$ws.Range('E23').Value = '  -0.67%  '
$ws.Range('D24').Value = '263.10'
$ws.Range('E24').Value = '  -3.77%  '
$ws.Range('D25').Value = '2.71'
$ws.Range('E25').Value = '  -1.87%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = '25.80'
$ws.Range('E27').Value = '  -2.46%  '
$ws.Range('E28').Value = '  +13.28%  '
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('E30').Value = '  -1.37%  '
$ws.Range('D31').Value = '51.54'
$ws.Range('E31').Value = '  +0.70%  '
$ws.Range('D32').Value = '34.38'
$ws.Range('E32').Value = '  +1.03%  '
$ws.Range('D33').Value = '5.99'
$ws.Range('E33').Value = '  +4.45%  '
$ws.Range('E34').Value = '  -7.60%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = '5.22'
$ws.Range('E35').Value = '  -2.97%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '0.0825'
$ws.Range('E36').Value = '  -1.53%  '
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('D38').Value = '18.25'
$ws.Range('E38').Value = '  +1.09%  '
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('E40').Value = '  -2.94%  '
$ws.Range('E41').Value = '  -1.42%  '
$ws.Range('E42').Value = '  -2.64%  '
$ws.Range('D43').Value = '120.83'
$ws.Range('E43').Value = '  -3.35%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '21.97'
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = '2.19'
$ws.Range('E45').Value = '  -2.46%  '
$ws.Range('D46').Value = '2.087.34'
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('E48').Value = '  -2.01%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').Value = '5.42'
$ws.Range('E49').Value = '  -4.83%  '
$ws.Range('B50').Value = 'SEI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range('D50').Value = '0.906'
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('E51').Value = '  +4.67%  '
